$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-CellText 'D2' '95.963.66'
Set-CellText 'E2' '  -1.95%  '
Set-CellText 'D3' '3.313.88'
Set-CellText 'E3' '  -2.96%  '
Set-CellText 'E4' '  -0.10%  '
Set-CellText 'D5' '247.87'
Set-CellText 'E5' '  -3.20%  '
Set-CellText 'D6' '648.98'
Set-CellText 'E6' '  -1.07%  '
Set-CellText 'E7' '  -9.62%  '
Set-CellText 'E8' '  -3.61%  '
Set-CellText 'E9' '  +0.09%  '
Set-CellText 'D10' '0.969'
Set-CellText 'E10' '  -8.00%  '
Set-CellText 'D11' '3.311.11'
Set-CellText 'E11' '  -2.90%  '
Set-CellText 'E12' '  -4.35%  '
Set-CellText 'D13' '39.58'
Set-CellText 'E13' '  -5.39%  '
Set-CellText 'D14' '95.796.38'
Set-CellText 'E14' '  -1.82%  '
Set-CellText 'E15' '  -5.29%  '
Set-CellText 'E16' '  -4.01%  '
Set-CellText 'D17' '3.932.00'
Set-CellText 'E17' '  -2.76%  '
Set-CellText 'D18' '8.40'
Set-CellText 'E18' '  -1.54%  '
Set-CellText 'D19' '3.310.75'
Set-CellText 'E19' '  -3.07%  '
Set-CellText 'D20' '16.77'
Set-CellText 'E20' '  -4.81%  '
Set-CellText 'D21' '0.478'
Set-CellText 'E21' '  -6.55%  '
Set-CellText 'D22' '500.06'
Set-CellText 'E22' '  -1.87%  '
Set-CellText 'D23' '10.37'
Set-CellText 'E23' '  -6.15%  '
Set-CellText 'E24' '  -4.67%  '
Set-CellText 'D25' '0.0000195'
Set-CellText 'E25' '  -5.60%  '
Set-CellText 'D26' '6.41'
Set-CellText 'E26' '  +5.36%  '
Set-CellText 'D27' '94.35'
Set-CellText 'E27' '  -1.88%  '
Set-CellText 'E28' '  -7.09%  '
Set-CellText 'D29' '3.487.06'
Set-CellText 'E29' '  -3.31%  '
Set-CellText 'E31' '  -9.81%  '
Set-CellText 'D32' '10.94'
Set-CellText 'E32' '  -4.81%  '
Set-CellText 'E33' '  -7.42%  '
Set-CellText 'E34' '  +9.19%  '
Set-CellText 'E35' '  -0.13%  '
Set-CellText 'D36' '0.539'
Set-CellText 'E36' '  -5.83%  '
Set-CellText 'D37' '27.68'
Set-CellText 'E37' '  -7.28%  '
Set-CellText 'E38' '  +4.48%  '
Set-CellText 'D39' '7.44'
Set-CellText 'E39' '  -4.55%  '
Set-CellText 'E40' '  +0.02%  '
Set-CellText 'D41' '0.149'
Set-CellText 'E41' '  -4.57%  '
Set-CellText 'D42' '500.93'
Set-CellText 'E42' '  -2.57%  '
Set-CellText 'E43' '  -1.88%  '
Set-CellText 'D44' '3.64'
Set-CellText 'E44' '  -1.91%  '
Set-CellText 'D45' '0.818'
Set-CellText 'E45' '  -4.81%  '
Set-CellText 'D46' '0.0405'
Set-CellText 'E46' '  -3.27%  '
Set-CellText 'D47' '5.47'
Set-CellText 'E47' '  -0.85%  '
Set-CellText 'B48' 'ImmutableX'
Set-CellText 'C48' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText 'D48' '1.62'
Set-CellText 'E48' '  +2.03%  '
Set-CellText 'B49' 'Cosmos'
Set-CellText 'C49' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText 'D49' '8.23'
Set-CellText 'D50' '52.81'
Set-CellText 'E50' '  +3.99%  '
Set-CellText 'E51' '  -6.42%  '
